$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2.257119139371683

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 0.8998700276285728

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 2.257119139371683
$ws.Range("D5").Value = 0.8998700276285728
